# Updated cryptos list on Tue May  7 03:41:19 UTC 2024 with GitHub Actions
#
# All data cells on Sheet1 are plain text (coin name / link / price / 1h
# volume-change), even the ones that look numeric, e.g. "5.77" or "0.535".
# Excel's normal COM `Range.Value =` setter auto-detects plain decimal-looking
# strings and silently turns them into Number cells, which would change the
# cell's stored type away from the Text type used throughout this sheet.
# `Set-TextValue` guards against that: it forces the cell to Text format,
# assigns the value, then clears the number format back off again so the
# cell keeps its original (unstyled) look, while the stored value stays Text.

$ws = $excel.ActiveWorkbook.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "63.126.36"
$ws.Range("E2").Value = "  -1.51%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.048.94"
$ws.Range("E3").Value = "  -3.15%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.17%  "

# Row 5 - BNB
Set-TextValue "D5" "586.17"
$ws.Range("E5").Value = "  -0.74%  "

# Row 6 - Solana
Set-TextValue "D6" "153.75"
$ws.Range("E6").Value = "  +5.06%  "

# Row 7
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
Set-TextValue "D8" "0.537"
$ws.Range("E8").Value = "  +1.23%  "

# Row 9
$ws.Range("D9").Value = "3.046.29"
$ws.Range("E9").Value = "  -2.89%  "

# Row 10
$ws.Range("E10").Value = "  -4.85%  "

# Row 11
Set-TextValue "D11" "5.78"
$ws.Range("E11").Value = "  -2.10%  "

# Row 12
Set-TextValue "D12" "0.449"
$ws.Range("E12").Value = "  -1.61%  "

# Row 13
Set-TextValue "D13" "36.68"
$ws.Range("E13").Value = "  -1.32%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  -4.86%  "

# Row 15 - was WrappedliquidstakedEther2.0, now TRON
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D15" "0.118"
$ws.Range("E15").Value = "  -2.01%  "

# Row 16 - was TRON, now WrappedliquidstakedEther2.0
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "3.555.79"
$ws.Range("E16").Value = "  -2.99%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "63.245.50"
$ws.Range("E17").Value = "  -0.94%  "

# Row 18 - was WrappedEther, now Polkadot
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D18" "7.09"
$ws.Range("E18").Value = "  -2.33%  "

# Row 19 - was Polkadot, now WrappedEther
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.042.49"
$ws.Range("E19").Value = "  -3.03%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "471.28"
$ws.Range("E20").Value = "  +1.18%  "

# Row 21
Set-TextValue "D21" "14.32"
$ws.Range("E21").Value = "  -0.50%  "

# Row 22
Set-TextValue "D22" "0.703"
$ws.Range("E22").Value = "  -3.95%  "

# Row 23
Set-TextValue "D23" "7.45"
$ws.Range("E23").Value = "  -1.82%  "

# Row 24
Set-TextValue "D24" "2.38"
$ws.Range("E24").Value = "  -0.52%  "

# Row 25
Set-TextValue "D25" "80.32"
$ws.Range("E25").Value = "  -0.72%  "

# Row 26
Set-TextValue "D26" "12.73"
$ws.Range("E26").Value = "  -3.33%  "

# Row 27
$ws.Range("E27").Value = "  +3.42%  "

# Row 28
Set-TextValue "D28" "0.999"
$ws.Range("E28").Value = "  -0.17%  "

# Row 29
Set-TextValue "D29" "7.46"
$ws.Range("E29").Value = "  +1.73%  "

# Row 30
Set-TextValue "D30" "0.997"
$ws.Range("E30").Value = "  -0.54%  "

# Row 31
Set-TextValue "D31" "2.64"
$ws.Range("E31").Value = "  -2.69%  "

# Row 32
Set-TextValue "D32" "2.14"
$ws.Range("E32").Value = "  -3.10%  "

# Row 33
$ws.Range("E33").Value = "  -1.75%  "

# Row 34
Set-TextValue "D34" "26.95"
$ws.Range("E34").Value = "  -2.55%  "

# Row 35
$ws.Range("D35").Value = "0.0₃0818"
$ws.Range("E35").Value = "  -4.74%  "

# Row 36
$ws.Range("E36").Value = "  -2.09%  "

# Row 37
Set-TextValue "D37" "3.31"
$ws.Range("E37").Value = "  +2.16%  "

# Row 38
Set-TextValue "D38" "5.95"
$ws.Range("E38").Value = "  -3.42%  "

# Row 39
Set-TextValue "D39" "2.18"
$ws.Range("E39").Value = "  -4.40%  "

# Row 40
Set-TextValue "D40" "9.28"
$ws.Range("E40").Value = "  -0.71%  "

# Row 41
Set-TextValue "D41" "50.47"
$ws.Range("E41").Value = "  -1.74%  "

# Row 42 - Bittensor (only price changes, % stays -6.08%)
Set-TextValue "D42" "435.08"

# Row 43
Set-TextValue "D43" "0.284"
$ws.Range("E43").Value = "  -2.68%  "

# Row 44
Set-TextValue "D44" "40.86"
$ws.Range("E44").Value = "  +2.33%  "

# Row 45
Set-TextValue "D45" "0.112"
$ws.Range("E45").Value = "  +3.54%  "

# Row 46
Set-TextValue "D46" "0.0357"
$ws.Range("E46").Value = "  -4.35%  "

# Row 47
$ws.Range("D47").Value = "2.781.90"
$ws.Range("E47").Value = "  -3.83%  "

# Row 48
Set-TextValue "D48" "130.27"
$ws.Range("E48").Value = "  -1.99%  "

# Row 50 (row 49, USDe, is untouched)
Set-TextValue "D50" "24.94"
$ws.Range("E50").Value = "  +3.70%  "

# Row 51
$ws.Range("E51").Value = "  -0.73%  "
